$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target Filename (D) and Impact (E) values for rows 4 through 82,
# reflecting the updated translation evaluation results.
$values = @(
    @("codeforces_96_B.py", "Runtime Error"),
    @("codeforces_306_A.py", "Runtime Error"),
    @("atcoder_ABC125_A.py", "Runtime Error"),
    @("atcoder_ABC051_A.py", "Runtime Error"),
    @("atcoder_ABC108_B.py", "Test Failed"),
    @("atcoder_ABC122_D.py", "Test Failed"),
    @("atcoder_ABC174_C.py", "Test Failed"),
    @("codeforces_99_A.py", "Test Failed"),
    @("codeforces_334_A.py", "Test Failed"),
    @("codeforces_49_A.py", "Test Failed"),
    @("codeforces_379_A.py", "Test Failed"),
    @("atcoder_ABC168_C.py", "Test Failed"),
    @("atcoder_ABC114_C.py", "Test Failed"),
    @("codeforces_340_A.py", "Test Failed"),
    @("atcoder_AGC046_A.py", "Test Failed"),
    @("atcoder_ABC172_D.py", "Test Failed"),
    @("codeforces_189_A.py", "Test Failed"),
    @("codeforces_190_A.py", "Test Failed"),
    @("atcoder_ABC139_B.py", "Test Failed"),
    @("atcoder_ABC142_A.py", "Test Failed"),
    @("codeforces_459_A.py", "Test Failed"),
    @("codeforces_544_B.py", "Test Failed"),
    @("atcoder_ABC132_F.py", "Test Failed"),
    @("atcoder_ABC132_A.py", "Test Failed"),
    @("atcoder_ABC124_C.py", "Test Failed"),
    @("atcoder_ABC164_A.py", "Test Failed"),
    @("codeforces_242_A.py", "Test Failed"),
    @("codeforces_678_B.py", "Test Failed"),
    @("atcoder_ABC178_A.py", "Test Failed"),
    @("codeforces_147_A.py", "Test Failed"),
    @("atcoder_ABC043_B.py", "Test Failed"),
    @("codeforces_369_B.py", "Test Failed"),
    @("codeforces_203_A.py", "Test Failed"),
    @("atcoder_ABC143_A.py", "Test Failed"),
    @("atcoder_AGC046_B.py", "Test Failed"),
    @("atcoder_ABC169_D.py", "Test Failed"),
    @("atcoder_ABC151_A.py", "Test Failed"),
    @("codeforces_79_A.py", "Test Failed"),
    @("codeforces_55_A.py", "Test Failed"),
    @("atcoder_ABC070_B.py", "Test Failed"),
    @("codeforces_30_A.py", "Test Failed"),
    @("codeforces_59_A.py", "Test Failed"),
    @("atcoder_ABC158_B.py", "Test Failed"),
    @("codeforces_622_A.py", "Test Failed"),
    @("codeforces_678_A.py", "Test Failed"),
    @("codeforces_276_B.py", "Test Failed"),
    @("atcoder_AGC002_A.py", "Test Failed"),
    @("atcoder_ABC127_B.py", "Test Failed"),
    @("atcoder_ABC149_B.py", "Test Failed"),
    @("codeforces_546_A.py", "Test Failed"),
    @("codeforces_569_A.py", "Test Failed"),
    @("atcoder_ABC120_C.py", "Test Failed"),
    @("atcoder_ABC149_C.py", "Test Failed"),
    @("codeforces_373_B.py", "Test Failed"),
    @("atcoder_ABC042_A.py", "Test Failed"),
    @("codeforces_86_A.py", "Test Failed"),
    @("codeforces_581_A.py", "Test Failed"),
    @("codeforces_579_A.py", "Test Failed"),
    @("codeforces_110_B.py", "Test Failed"),
    @("atcoder_ABC155_E.py", "Test Failed"),
    @("atcoder_ABC124_A.py", "Test Failed"),
    @("atcoder_ABC170_A.py", "Test Failed"),
    @("atcoder_AGC025_A.py", "Test Failed"),
    @("atcoder_ABC169_C.py", "Test Failed"),
    @("codeforces_8_B.py", "Test Failed"),
    @("atcoder_ABC158_A.py", "Test Failed"),
    @("codeforces_171_A.py", "Test Failed"),
    @("atcoder_ABC136_B.py", "Test Failed"),
    @("codeforces_92_A.py", "Test Failed"),
    @("codeforces_514_A.py", "Test Failed"),
    @("codeforces_32_B.py", "Test Failed"),
    @("atcoder_ABC153_A.py", "Test Failed"),
    @("codeforces_58_B.py", "Test Failed"),
    @("codeforces_672_A.py", "Test Failed"),
    @("codeforces_651_A.py", "Test Failed"),
    @("atcoder_ARC062_B.py", "Test Failed"),
    @("atcoder_ARC102_C.py", "Infinite Loop"),
    @("codeforces_669_A.py", "Infinite Loop"),
    @("atcoder_AGC006_B.py", "Correct")
)

$startRow = 4
for ($i = 0; $i -lt $values.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 4).Value = $values[$i][0]
    $ws.Cells.Item($row, 5).Value = $values[$i][1]
}

$wb.Save()
